# Updated symbol list on Wed Dec 28 23:36:06 UTC 2022 with GitHub Actions
# Applies refreshed price/volume data scraped for the cryptos sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as TEXT (preserving the original inline-string /
# shared-string cell type) instead of letting Excel auto-convert
# numeric-looking text into a real number. We do this by prefixing the
# value with a leading apostrophe (Excel's "treat as text" marker) and
# then restoring the cell's original Style, since applying the value can
# implicitly attach a "quote prefix" style to the cell.
function Set-TextValue {
    param($Sheet, [string]$Address, [string]$Val)
    $cell = $Sheet.Range($Address)
    $origStyle = $cell.Style
    $cell.Value = "'" + $Val
    $cell.Style = $origStyle
}

# --- Price column (D) updates ---
Set-TextValue $ws "D2"  "243.49"
Set-TextValue $ws "D3"  "23.79"
Set-TextValue $ws "D4"  "5.236"
Set-TextValue $ws "D5"  "0.05760"
Set-TextValue $ws "D6"  "6.405"
Set-TextValue $ws "D7"  "3.243"
Set-TextValue $ws "D8"  "0.8126"
Set-TextValue $ws "D9"  "0.8806"
Set-TextValue $ws "D11" "0.07064"
Set-TextValue $ws "D12" "0.03149"
Set-TextValue $ws "D13" "0.03042"
Set-TextValue $ws "D14" "0.09315"
Set-TextValue $ws "D15" "3.805"
Set-TextValue $ws "D16" "0.001529"
Set-TextValue $ws "D17" "0.04696"
Set-TextValue $ws "D18" "0.0006013"
Set-TextValue $ws "D19" "0.006201"
Set-TextValue $ws "D20" "0.001236"
Set-TextValue $ws "D21" "0.004068"
Set-TextValue $ws "D22" "0.00008689"
Set-TextValue $ws "D23" "3.545"
Set-TextValue $ws "D24" "2.145"
Set-TextValue $ws "D40" "0.03730"
Set-TextValue $ws "D41" "0.006242"
Set-TextValue $ws "D42" "0.1045"
Set-TextValue $ws "D43" "0.002597"
Set-TextValue $ws "D44" "0.007954"
Set-TextValue $ws "D45" "0.00005292"
Set-TextValue $ws "D47" "0.5294"
Set-TextValue $ws "D48" "0.002505"
Set-TextValue $ws "D49" "0.00002097"

# --- Volume(1h) column (E) updates (Best/Worst-in-24h badge toggles) ---
Set-TextValue $ws "E18" "17OneONE"
Set-TextValue $ws "E41" "40KickTokenKICKBestin24h"
Set-TextValue $ws "E47" "46CoinbaseStockTokenCOINWorstin24h"
Set-TextValue $ws "E48" "47BOLOBOLO"
